# Update "want to go" counts (column F) across the workbook's sheets
# to reflect the latest scrape snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 4340
$wsExpo.Range("F6").Value = 420
$wsExpo.Range("F7").Value = 3450
$wsExpo.Range("F8").Value = 978
$wsExpo.Range("F11").Value = 298
$wsExpo.Range("F12").Value = 2340
$wsExpo.Range("F13").Value = 1252
$wsExpo.Range("F16").Value = 500
$wsExpo.Range("F17").Value = 245
$wsExpo.Range("F18").Value = 52
$wsExpo.Range("F19").Value = 9716
$wsExpo.Range("F20").Value = 5923
$wsExpo.Range("F21").Value = 377
$wsExpo.Range("F23").Value = 810
$wsExpo.Range("F24").Value = 118
$wsExpo.Range("F25").Value = 825
$wsExpo.Range("F26").Value = 3510
$wsExpo.Range("F29").Value = 449
$wsExpo.Range("F30").Value = 103
$wsExpo.Range("F31").Value = 229
$wsExpo.Range("F32").Value = 210
$wsExpo.Range("F33").Value = 4777
$wsExpo.Range("F34").Value = 16
$wsExpo.Range("F35").Value = 1032
$wsExpo.Range("F36").Value = 133
$wsExpo.Range("F37").Value = 11
$wsExpo.Range("F38").Value = 455

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 8641
$wsLocal.Range("F4").Value = 1515

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8641
$wsAll.Range("F5").Value = 1515
$wsAll.Range("F7").Value = 4340
$wsAll.Range("F9").Value = 420
$wsAll.Range("F10").Value = 3450
$wsAll.Range("F11").Value = 978
$wsAll.Range("F14").Value = 298
$wsAll.Range("F15").Value = 2340
$wsAll.Range("F19").Value = 1252
$wsAll.Range("F23").Value = 500
$wsAll.Range("F24").Value = 245
$wsAll.Range("F25").Value = 52
$wsAll.Range("F26").Value = 9716
$wsAll.Range("F29").Value = 377
$wsAll.Range("F31").Value = 810
$wsAll.Range("F32").Value = 118
$wsAll.Range("F33").Value = 825
$wsAll.Range("F34").Value = 3510
$wsAll.Range("F37").Value = 449
$wsAll.Range("F38").Value = 103
$wsAll.Range("F39").Value = 229
$wsAll.Range("F41").Value = 210
$wsAll.Range("F42").Value = 4777
$wsAll.Range("F43").Value = 1032
$wsAll.Range("F44").Value = 133
$wsAll.Range("F45").Value = 455
